$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "PROTEOMICS_EOTAXIN (HUMAN)"
$ws.Range("B15").Value = "Eotaxin, human"

$ws.Range("A16").Value = "PROTEOMICS_IFN-Y (HUMAN)"
$ws.Range("B16").Value = "Interferon y, human"

$ws.Range("A17").Value = "PROTEOMICS_IL-15"
$ws.Range("B17").Value = "Interleukin-15"

$ws.Range("A18").Value = "PROTEOMICS_MCP-4 (HUMAN)"
$ws.Range("B18").Value = "Monocyte chemotactic protein-4, human"

$ws.Range("A19").Value = "PROTEOMICS_MDC (HUMAN)"
$ws.Range("B19").Value = "Myeloid dentritic cells, human"

$ws.Range("A20").Value = "PROTEOMICS_MIP-1A (HUMAN)"
$ws.Range("B20").Value = "Macrophage Inflammatory Protein 1a, human"

$ws.Range("A21").Value = "PROTEOMICS_SAA"
$ws.Range("B21").Value = "Serum amyloid A"

$ws.Range("A22").Value = "PROTEOMICS_VEGF- ANGIO PLATE"
$ws.Range("B22").Value = "Vascular endothelial growth factor angio plate"
